$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 2.38
$ws.Range("AB2").Value = 1.53
$ws.Range("AD2").Value = 6.5
$ws.Range("AF2").Value = 13
$ws.Range("AN2").Value = 11
$ws.Range("AO2").Value = 26
$ws.Range("AP2").Value = 19
$ws.Range("AQ2").Value = 67
$ws.Range("AR2").Value = 51
$ws.Range("AS2").Value = 67
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.5
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("U2").Value = 4.1
$ws.Range("V2").Value = 1.22
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("AA4").Value = 2.38
$ws.Range("AB4").Value = 1.53
$ws.Range("AC4").Value = 6.5
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 34
$ws.Range("AG4").Value = 34
$ws.Range("AH4").Value = 51
$ws.Range("AI4").Value = 5
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 81
$ws.Range("AN4").Value = 6
$ws.Range("AO4").Value = 11
$ws.Range("AP4").Value = 12
$ws.Range("AR4").Value = 29
$ws.Range("AS4").Value = 51
$ws.Range("G4").Value = 2.9
$ws.Range("H4").Value = 2.8
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 1.8
$ws.Range("L4").Value = 3.75
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.78
$ws.Range("S4").Value = 3.1
$ws.Range("T4").Value = 1.36
$ws.Range("U4").Value = 5.2
$ws.Range("V4").Value = 1.15
$ws.Range("W4").Value = 6.5
$ws.Range("X4").Value = 1.11
$ws.Range("Y4").Value = 1.67
$ws.Range("Z4").Value = 2.1
$ws.Range("AC5").Value = 6.1
$ws.Range("AH5").Value = 65
$ws.Range("AI5").Value = 4.2
$ws.Range("AN5").Value = 5.6
$ws.Range("AP5").Value = 11.75
$ws.Range("AR5").Value = 37
$ws.Range("AS5").Value = 65
$ws.Range("J5").Value = 3.9
$ws.Range("K5").Value = 1.72
$ws.Range("L5").Value = 3.7
$ws.Range("N5").Value = 4.2
$ws.Range("Y5").Value = 1.72
$ws.Range("Z5").Value = 2
$ws.Range("AA6").Value = 2.3
$ws.Range("AB6").Value = 1.55
$ws.Range("AC6").Value = 6.2
$ws.Range("AE6").Value = 11.5
$ws.Range("AG6").Value = 35
$ws.Range("AH6").Value = 60
$ws.Range("AI6").Value = 4.55
$ws.Range("AK6").Value = 20
$ws.Range("AP6").Value = 11.75
$ws.Range("AS6").Value = 65
$ws.Range("I6").Value = 2.77
$ws.Range("J6").Value = 3.6
$ws.Range("K6").Value = 1.83
$ws.Range("L6").Value = 3.55
$ws.Range("M6").Value = 1.16
$ws.Range("N6").Value = 4.55
$ws.Range("O6").Value = 1.7
$ws.Range("P6").Value = 2.05
$ws.Range("S6").Value = 2.95
$ws.Range("T6").Value = 1.35
$ws.Range("W6").Value = 5.5
$ws.Range("X6").Value = 1.11
$ws.Range("Y6").Value = 1.62
$ws.Range("Z6").Value = 2.18
$ws.Range("W7").Value = 6.5
$ws.Range("X7").Value = 1.11
$ws.Range("AA9").Value = 2.05
$ws.Range("AB9").Value = 1.7
$ws.Range("AC9").Value = 6.5
$ws.Range("AD9").Value = 10
$ws.Range("AE9").Value = 10
$ws.Range("AG9").Value = 23
$ws.Range("AH9").Value = 41
$ws.Range("AI9").Value = 7.5
$ws.Range("AK9").Value = 19
$ws.Range("AL9").Value = 67
$ws.Range("AM9").Value = 501
$ws.Range("AN9").Value = 7.5
$ws.Range("AQ9").Value = 34
$ws.Range("AR9").Value = 29
$ws.Range("AS9").Value = 41
$ws.Range("G9").Value = 2.38
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 2.9
$ws.Range("J9").Value = 3.25
$ws.Range("K9").Value = 1.95
$ws.Range("L9").Value = 3.75
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 1.44
$ws.Range("P9").Value = 2.63
$ws.Range("Q9").Value = 1.85
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2.4
$ws.Range("T9").Value = 1.53
$ws.Range("W9").Value = 5
$ws.Range("X9").Value = 1.17
$ws.Range("Y9").Value = 1.53
$ws.Range("Z9").Value = 2.38
$ws.Range("Q10").Value = 1.98
$ws.Range("R10").Value = 1.88
$ws.Range("S10").Value = 2.6
$ws.Range("T10").Value = 1.48
$ws.Range("W10").Value = 5.5
$ws.Range("X10").Value = 1.14
$ws.Range("AI11").Value = 7.5
$ws.Range("N11").Value = 7.5
$ws.Range("O11").Value = 1.4
$ws.Range("P11").Value = 2.75
$ws.Range("AA12").Value = 2.05
$ws.Range("AB12").Value = 1.7
$ws.Range("AC12").Value = 6
$ws.Range("AD12").Value = 9.5
$ws.Range("AE12").Value = 10
$ws.Range("AG12").Value = 21
$ws.Range("AH12").Value = 41
$ws.Range("AI12").Value = 7
$ws.Range("AJ12").Value = 6
$ws.Range("AK12").Value = 19
$ws.Range("AL12").Value = 67
$ws.Range("AM12").Value = 501
$ws.Range("AN12").Value = 7.5
$ws.Range("AR12").Value = 34
$ws.Range("G12").Value = 2.25
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 3.25
$ws.Range("J12").Value = 3.1
$ws.Range("K12").Value = 1.95
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.5
$ws.Range("Q12").Value = 1.9
$ws.Range("R12").Value = 1.95
$ws.Range("S12").Value = 2.5
$ws.Range("T12").Value = 1.5
$ws.Range("W12").Value = 5
$ws.Range("X12").Value = 1.17
$ws.Range("Y12").Value = 1.57
$ws.Range("Z12").Value = 2.25
$ws.Range("AA13").Value = 2.05
$ws.Range("AB13").Value = 1.7
$ws.Range("AC13").Value = 6.5
$ws.Range("AM13").Value = 401
$ws.Range("AQ13").Value = 101
$ws.Range("I13").Value = 8
$ws.Range("J13").Value = 1.95
$ws.Range("L13").Value = 7.5
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("O13").Value = 1.25
$ws.Range("P13").Value = 3.75
$ws.Range("S13").Value = 1.88
$ws.Range("T13").Value = 1.98
$ws.Range("W13").Value = 3.25
$ws.Range("X13").Value = 1.33
$ws.Range("AC14").Value = 7.5
$ws.Range("AD14").Value = 10
$ws.Range("AF14").Value = 19
$ws.Range("G14").Value = 2.1
$ws.Range("I14").Value = 3.25
$ws.Range("J14").Value = 2.75
$ws.Range("L14").Value = 3.75
$ws.Range("S14").Value = 2.05
$ws.Range("T14").Value = 1.8
$ws.Range("W14").Value = 3.5
$ws.Range("X14").Value = 1.29
